# Fill in the May 2023 statistics values (Circulation / ILL Loans / ILL Borrows)
# for every library row (3-59) on the single worksheet. These cells were
# previously blank placeholders and are now populated with the reported numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, column B (Circulation), column C (ILL Loans), column D (ILL Borrows)
$data = @(
    @(3, 28933, 4523, 4864),
    @(4, 14555, 1840, 1730),
    @(5, 50595, 4591, 4652),
    @(6, 896, 431, 108),
    @(7, 31748, 5597, 4645),
    @(8, 3887, 872, 850),
    @(9, 3901, 784, 479),
    @(10, 1716, 287, 191),
    @(11, 119, 143, 29),
    @(12, 3, 0, 0),
    @(13, 608, 143, 211),
    @(14, 2037, 859, 617),
    @(15, 3519, 1342, 567),
    @(16, 2186, 839, 322),
    @(17, 1362, 578, 123),
    @(18, 11371, 1932, 2173),
    @(19, 1219, 437, 382),
    @(20, 12633, 1550, 2241),
    @(21, 133, 256, 10),
    @(22, 11407, 1571, 2168),
    @(23, 813, 260, 104),
    @(24, 12609, 1986, 2647),
    @(25, 54315, 5059, 6688),
    @(26, 3676, 1281, 711),
    @(27, 0, 0, 0),
    @(28, 3502, 758, 903),
    @(29, 954, 311, 195),
    @(30, 9747, 1893, 1653),
    @(31, 318, 140, 160),
    @(32, 1678, 1181, 298),
    @(33, 11204, 2403, 1876),
    @(34, 6477, 2080, 1540),
    @(35, 4018, 454, 1031),
    @(36, 36525, 4013, 3794),
    @(37, 5465, 1993, 862),
    @(38, 16409, 1450, 1956),
    @(39, 551, 636, 163),
    @(40, 1391, 337, 479),
    @(41, 2228, 271, 93),
    @(42, 8333, 481, 265),
    @(43, 236, 86, 66),
    @(44, 536, 38, 43),
    @(45, 1045, 14, 7),
    @(46, 2192, 642, 273),
    @(47, 7721, 2394, 1425),
    @(48, 21835, 2486, 3268),
    @(49, 9317, 2530, 843),
    @(50, 7684, 759, 1101),
    @(51, 21009, 2100, 3152),
    @(52, 3326, 396, 927),
    @(53, 10212, 2125, 1772),
    @(54, 1093, 831, 465),
    @(55, 1463, 952, 112),
    @(56, 2090, 589, 644),
    @(57, 8350, 3183, 1717),
    @(58, 12378, 890, 474),
    @(59, 452800, 71504, 65878),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

